$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new logged entry as row 23 (Start Time / Finish Time / Minutes Logged)
$ws.Range("A23").Value = "12:08PM 1-20-2018"
$ws.Range("B23").Value = "2:21PM 1-20-2018"
$ws.Range("C23").Value = 133

# Update the selected cell to reflect where the user was working (C26)
$ws.Range("C26").Select() | Out-Null
